$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values in row 1 (columns B-E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row (row 2) values for columns B-E
$ws.Range("B2").Value = 0.71303677533663179
$ws.Range("C2").Value = 1.9134581246704694
$ws.Range("D2").Value = 1.0810840941366469
$ws.Range("E2").Value = 2.0743975788486257

# Update STR row (row 3) values for columns B-E
$ws.Range("B3").Value = 1.9833558962570397
$ws.Range("C3").Value = 5.803810009943482
$ws.Range("D3").Value = 5.0498057792531421
$ws.Range("E3").Value = 2.571109060183042

# Update the selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
